# Fix/populate the course-info sheet for LOM3212 to match the published
# syllabus data: corrects label/content misalignment in rows 10-23 and
# appends the new "Docentes responsaveis", "Programa"/"Programa resumido"
# bodies, "Bibliografia" text and the "Requisitos" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A was "min=1 max=2" (redundant with the col-2 override that
# follows it) - narrow it down to just column 1, matching column 2's
# own width/style definition.
$ws.Columns.Item(1).ColumnWidth = 29.83

# Row 1

# Row 2

# Row 3

# Row 4

# Row 5

# Row 6

# Row 7

# Row 8

# Row 9

# Row 10
$ws.Range("B10").Value = 'Apresentar noções de mecânica dos fluidos, mediante estudo dos meios fluidos quando estáticos ou em movimento. Capacitar o aluno a modelar e resolver problemas de interesse em mecânica dos fluidos, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução.'
$ws.Range("C10").Value = 'Apresentar noções de mecânica dos fluidos, mediante estudo dos meios fluidos quando estáticos ou em movimento. Capacitar o aluno a modelar e resolver problemas de interesse em mecânica dos fluidos, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução.'

# Row 11

# Row 12

# Row 13
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Range("C13").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range("C14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Rows.Item(14).AutoFit()

# Row 15
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = 'Fundamentos de mecânica dos fluidos. Introdução à estática dos fluidos. Formulação integral e diferencial das equações de transporte de massa, energia e quantidade de movimento. Análise dimensional e semelhança. Escoamento incompressível de fluidos ideais e viscosos, regime laminar e turbulento. Equação de Navier-Stokes. Teoria da camada limite.'
$ws.Range("C15").Value = 'Fundamentos de mecânica dos fluidos. Introdução à estática dos fluidos. Formulação integral e diferencial das equações de transporte de massa, energia e quantidade de movimento. Análise dimensional e semelhança. Escoamento incompressível de fluidos ideais e viscosos, regime laminar e turbulento. Equação de Navier-Stokes. Teoria da camada limite.'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = 'Fundamentals of fluid mechanics. Introduction to fluid statics. Integral and differential formulation of mass, energy and momentum transport equations. Dimensional analysis and similarity. Incompressible flow of ideal and viscous fluids, laminar and turbulent regime. Navier-Stokes equation. Boundary layer theory.'
$ws.Range("C16").Value = 'Fundamentals of fluid mechanics. Introduction to fluid statics. Integral and differential formulation of mass, energy and momentum transport equations. Dimensional analysis and similarity. Incompressible flow of ideal and viscous fluids, laminar and turbulent regime. Navier-Stokes equation. Boundary layer theory.'
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = 'Introdução: conceito de fluido; propriedades e conceito de contínuo; modelagem de processos de transferência; métodos de análise; dimensões e unidades.Revisão de estática de fluidos: equação básica da hidrostática, variação de pressão em um fluido estático; princípios de Stevin, de Pascal e de Arquimedes.Formulação integral das equações de transporte: teorema de transporte de Reynolds; aplicação para os princípios de conservação de massa, quantidade de movimento e energia; equação de Bernoulli.Formulação diferencial das equações de transporte: descrição do escoamento; forma diferencial: dos princípios de conservação de massa, quantidade de movimento e energia; formulação adimensional, análise dimensional e semelhança. Grupos adimensionais: número de Reynolds e número de Grashoff.Escoamento incompressível interno: equações de Euler; lei de Newton para a viscosidade, tensões de cisalhamento; equação de Navier-Stokes; regimes de escoamento: escoamento laminar e turbulento. Cálculo de perda de carga (distribuída e localizada), coeficiente de atrito. Escoamento incompressível externo: introdução à camada limite; escoamento ao redor de corpos, força da arraste.'
$ws.Range("B17").Font.Bold = $false
$ws.Range("B17").WrapText = $true
$ws.Range("B17").VerticalAlignment = -4160
$ws.Range("C17").Value = 'Introdução: conceito de fluido; propriedades e conceito de contínuo; modelagem de processos de transferência; métodos de análise; dimensões e unidades.Revisão de estática de fluidos: equação básica da hidrostática, variação de pressão em um fluido estático; princípios de Stevin, de Pascal e de Arquimedes.Formulação integral das equações de transporte: teorema de transporte de Reynolds; aplicação para os princípios de conservação de massa, quantidade de movimento e energia; equação de Bernoulli.Formulação diferencial das equações de transporte: descrição do escoamento; forma diferencial: dos princípios de conservação de massa, quantidade de movimento e energia; formulação adimensional, análise dimensional e semelhança. Grupos adimensionais: número de Reynolds e número de Grashoff.Escoamento incompressível interno: equações de Euler; lei de Newton para a viscosidade, tensões de cisalhamento; equação de Navier-Stokes; regimes de escoamento: escoamento laminar e turbulento. Cálculo de perda de carga (distribuída e localizada), coeficiente de atrito. Escoamento incompressível externo: introdução à camada limite; escoamento ao redor de corpos, força da arraste.'
$ws.Range("C17").Font.Bold = $false
$ws.Range("C17").WrapText = $true
$ws.Range("C17").VerticalAlignment = -4160
$ws.Range("C17").Font.Color = 255
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = 'Introduction: fluid concept; properties and concept of continuum; modeling of transfer processes; analysis methods; dimensions and units. • Fluid statics review: basic hydrostatic equation, pressure variation in a static fluid; principles of Stevin, Pascal and Archimedes. • Integral formulation of transport equations: Reynolds transport theorem; application to the principles of conservation of mass, momentum and energy; Bernoulli equation. • Differential formulation of transport equations: description of the flow; differential form: from the principles of conservation of mass, momentum and energy; dimensionless formulation, dimensional analysis and similarity. Dimensionless groups: Reynolds number and Grashoff number. • Internal incompressible flow: Euler equations; Newton''s law for viscosity, shear stresses; Navier-Stokes equation; flow regimes: laminar and turbulent flow. Calculation of pressure drop (distributed and localized), friction coefficient. • External incompressible flow: introduction to the boundary layer; flow around bodies, drag force.'
$ws.Range("C18").Value = 'Introduction: fluid concept; properties and concept of continuum; modeling of transfer processes; analysis methods; dimensions and units. • Fluid statics review: basic hydrostatic equation, pressure variation in a static fluid; principles of Stevin, Pascal and Archimedes. • Integral formulation of transport equations: Reynolds transport theorem; application to the principles of conservation of mass, momentum and energy; Bernoulli equation. • Differential formulation of transport equations: description of the flow; differential form: from the principles of conservation of mass, momentum and energy; dimensionless formulation, dimensional analysis and similarity. Dimensionless groups: Reynolds number and Grashoff number. • Internal incompressible flow: Euler equations; Newton''s law for viscosity, shear stresses; Navier-Stokes equation; flow regimes: laminar and turbulent flow. Calculation of pressure drop (distributed and localized), friction coefficient. • External incompressible flow: introduction to the boundary layer; flow around bodies, drag force.'
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Rows.Item(19).AutoFit()

# Row 20
$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Range("C20").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Range("C21").Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("B22").Font.Bold = $false
$ws.Range("B22").WrapText = $true
$ws.Range("B22").VerticalAlignment = -4160
$ws.Range("C22").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C22").Font.Bold = $false
$ws.Range("C22").WrapText = $true
$ws.Range("C22").VerticalAlignment = -4160
$ws.Range("C22").Font.Color = 255
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").WrapText = $false
$ws.Range("A23").VerticalAlignment = -4160
$ws.Range("B23").Value = 'BIRD,R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. LTC Editora, 2004.
FOX, R. W., McDONALD, A. T. Introdução à Mecânica dos Fluidos. LTC Editora, 2001.
SISSOM, L. E., PITTS, D. R. Fenômenos de Transporte. Ed. Guanabara, 1988.'
$ws.Range("C23").Value = 'BIRD,R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. LTC Editora, 2004.
FOX, R. W., McDONALD, A. T. Introdução à Mecânica dos Fluidos. LTC Editora, 2001.
SISSOM, L. E., PITTS, D. R. Fenômenos de Transporte. Ed. Guanabara, 1988.'
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$ws.Range("A24").Value = 'Requisitos:'
$ws.Range("A24").Font.Bold = $true
$ws.Range("A24").WrapText = $false
$ws.Range("A24").VerticalAlignment = -4160

# Row 25
$ws.Range("B25").Value = 'LOB1019 -  Física II  (Requisito)
'
$ws.Range("B25").Font.Bold = $false
$ws.Range("B25").WrapText = $true
$ws.Range("B25").VerticalAlignment = -4160
$ws.Range("C25").Value = 'LOB1019 -  Física II  (Requisito)
'
$ws.Range("C25").Font.Bold = $false
$ws.Range("C25").WrapText = $true
$ws.Range("C25").VerticalAlignment = -4160
$ws.Range("C25").Font.Color = 255
$ws.Rows.Item(25).RowHeight = 30

